$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.585.68"
$ws.Range("E2").Value = "  -1.52%  "

$ws.Range("D3").Value = "3.407.45"
$ws.Range("E3").Value = "  -0.41%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "568.59"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").Value = "157.30"
$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "3.412.09"
$ws.Range("E8").Value = "  -0.41%  "

$ws.Range("D9").Value = "0.569"
$ws.Range("E9").Value = "  -7.97%  "

$ws.Range("E10").Value = "  +1.20%  "

$ws.Range("E11").Value = "  -4.25%  "

$ws.Range("D12").Value = "0.421"
$ws.Range("E12").Value = "  -4.57%  "

$ws.Range("D13").Value = "3.993.53"
$ws.Range("E13").Value = "  -0.45%  "

$ws.Range("E14").Value = "  +0.24%  "

$ws.Range("D15").Value = "26.91"
$ws.Range("E15").Value = "  -3.53%  "

$ws.Range("D16").Value = "0.0000171"
$ws.Range("E16").Value = "  -9.12%  "

$ws.Range("D17").Value = "63.639.01"
$ws.Range("E17").Value = "  -1.47%  "

$ws.Range("D18").Value = "3.374.60"
$ws.Range("E18").Value = "  -0.92%  "

$ws.Range("E19").Value = "  -4.68%  "

$ws.Range("E20").Value = "  -2.97%  "

$ws.Range("D21").Value = "385.35"
$ws.Range("E21").Value = "  +1.85%  "

$ws.Range("D22").Value = "7.74"
$ws.Range("E22").Value = "  -3.44%  "

$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Value = "71.14"
$ws.Range("E24").Value = "  -1.87%  "

$ws.Range("D25").Value = "0.514"
$ws.Range("E25").Value = "  -6.72%  "

$ws.Range("D26").Value = "0.0000114"
$ws.Range("E26").Value = "  -4.94%  "

$ws.Range("D27").Value = "9.65"
$ws.Range("E27").Value = "  -6.11%  "

$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").Value = "6.06"
$ws.Range("E30").Value = "  -2.72%  "

$ws.Range("D31").Value = "1.38"
$ws.Range("E31").Value = "  -6.74%  "

$ws.Range("E32").Value = "  -2.65%  "

$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("D34").Value = "22.83"
$ws.Range("E34").Value = "  -1.32%  "

$ws.Range("D35").Value = "6.93"
$ws.Range("E35").Value = "  -4.23%  "

$ws.Range("D36").Value = "1.50"
$ws.Range("E36").Value = "  -6.96%  "

$ws.Range("D37").Value = "160.70"
$ws.Range("E37").Value = "  +0.67%  "

$ws.Range("E38").Value = "  +8.70%  "

$ws.Range("E39").Value = "  -4.91%  "

$ws.Range("D40").Value = "2.785.80"
$ws.Range("E40").Value = "  -3.15%  "

$ws.Range("D41").Value = "25.83"
$ws.Range("E41").Value = "  -3.90%  "

$ws.Range("D42").Value = "0.0722"
$ws.Range("E42").Value = "  -5.57%  "

$ws.Range("D43").Value = "42.96"
$ws.Range("E43").Value = "  +0.37%  "

$ws.Range("E44").Value = "  -8.81%  "

$ws.Range("D45").Value = "25.77"
$ws.Range("E45").Value = "  -3.37%  "

$ws.Range("D46").Value = "4.33"
$ws.Range("E46").Value = "  -6.12%  "

$ws.Range("E47").Value = "  -4.71%  "

$ws.Range("D48").Value = "2.35"
$ws.Range("E48").Value = "  +7.96%  "

$ws.Range("D49").Value = "326.82"
$ws.Range("E49").Value = "  +2.27%  "

$ws.Range("E50").Value = "  -5.34%  "

$ws.Range("E51").Value = "  -4.96%  "
